$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the login credentials in rows 2-4 (row 5 keeps the old values)
# Write column B first so the shared-strings table gets "sEvYmEq" (index 4)
# before "mngr601962" (index 5), matching the target ordering.
$ws.Range("B2:B4").Value = "sEvYmEq"
$ws.Range("A2:A4").Value = "mngr601962"

# Update the selected range / active cell shown in the sheet view
$ws.Range("A2:B4").Select()
